$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 3 (the "DIETEL PARTNERS" record), shifting the rows
# below it (row 4, "OAK VIEW NATIONAL BANK") up by one.
$ws.Rows("3").Delete()

# Match the author's final cell selection recorded in the saved file.
$ws.Range("D7").Select()
